$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), O (Origen), P (Precio $/Kg)
# These reflect the weekly refresh of the dataset (rows reshuffled by date).

$rows = @(
    @{ Row=2;  D=44292; J=90;  K=6000; L=6000; M=6000; O="Región Metropolitana"; P=375 },
    @{ Row=3;  D=44232; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=4;  D=44251; J=120; K=5000; L=5000; M=5000; O="Región Metropolitana"; P=312 },
    @{ Row=5;  D=44186; J=160; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=6;  D=44204; J=430; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=7;  D=44208; J=160; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=8;  D=44210; J=340; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=9;  D=44189; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=10; D=44230; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=11; D=44215; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=12; D=44188; J=210; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=13; D=44187; J=160; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 },
    @{ Row=14; D=44231; J=250; K=5000; L=6000; M=5500; O="Provincia de Quillota"; P=344 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $r.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 15).Value = $r.O   # O - Origen
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio $/Kg
}
